$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 5509
$ws.Cells.Item(18, 9).Value = 3073.875
$ws.Cells.Item(18, 10).Value = 24990
$ws.Cells.Item(18, 11).Value = 3073.875
$ws.Cells.Item(18, 12).Value = 24990
$ws.Cells.Item(18, 13).Value = -2789.875
$ws.Cells.Item(18, 14).Value = -25558
$ws.Cells.Item(28, 8).Value = 1088.174
$ws.Cells.Item(28, 9).Value = 1001.5263
$ws.Cells.Item(28, 11).Value = 1001.5263
$ws.Cells.Item(28, 13).Value = -516.5263
$ws.Cells.Item(29, 8).Value = 2481.6
$ws.Cells.Item(29, 9).Value = 1555
$ws.Cells.Item(29, 11).Value = 4665
$ws.Cells.Item(29, 13).Value = -4384
$ws.Cells.Item(32, 8).Value = 7935
$ws.Cells.Item(32, 9).Value = 7784.6
$ws.Cells.Item(32, 11).Value = 7784.6
$ws.Cells.Item(32, 13).Value = -7458.6
$ws.Cells.Item(33, 8).Value = 1607.56
$ws.Cells.Item(33, 9).Value = 1823.1
$ws.Cells.Item(33, 11).Value = 1823.1
$ws.Cells.Item(33, 13).Value = -1594.1
$ws.Cells.Item(41, 8).Value = 1144.6666
$ws.Cells.Item(41, 9).Value = 1144.6666
$ws.Cells.Item(41, 11).Value = 1144.6666
$ws.Cells.Item(41, 13).Value = -704.6666
$ws.Cells.Item(43, 8).Value = 1184.625
$ws.Cells.Item(43, 9).Value = 998.1429
$ws.Cells.Item(43, 10).Value = 2490
$ws.Cells.Item(43, 11).Value = 998.1429
$ws.Cells.Item(43, 12).Value = 2490
$ws.Cells.Item(43, 13).Value = -929.1429
$ws.Cells.Item(43, 14).Value = -2628
$ws.Cells.Item(45, 8).Value = 5749.5
$ws.Cells.Item(45, 9).Value = 1500
$ws.Cells.Item(45, 11).Value = 4500
$ws.Cells.Item(45, 13).Value = -4308
$ws.Cells.Item(70, 8).Value = 947.0833
$ws.Cells.Item(70, 9).Value = 849
$ws.Cells.Item(70, 10).Value = 1437.5
$ws.Cells.Item(70, 11).Value = 2547
$ws.Cells.Item(70, 12).Value = 4312.5
$ws.Cells.Item(70, 13).Value = -2277
$ws.Cells.Item(70, 14).Value = -4852.5
$ws.Cells.Item(73, 8).Value = 947.0833
$ws.Cells.Item(73, 9).Value = 849
$ws.Cells.Item(73, 10).Value = 1437.5
$ws.Cells.Item(73, 11).Value = 2547
$ws.Cells.Item(73, 12).Value = 4312.5
$ws.Cells.Item(73, 13).Value = -1611
$ws.Cells.Item(73, 14).Value = -6184.5
$ws.Cells.Item(74, 8).Value = 12471.286
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(77, 8).Value = 12471.286
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(86, 8).Value = 2191.875
$ws.Cells.Item(86, 9).Value = 2510.6
$ws.Cells.Item(86, 11).Value = 2510.6
$ws.Cells.Item(86, 13).Value = -1387.6
$ws.Cells.Item(87, 8).Value = 59570.645
$ws.Cells.Item(87, 10).Value = 59570.645
$ws.Cells.Item(87, 12).Value = 59570.645
$ws.Cells.Item(87, 14).Value = -62066.645
$ws.Cells.Item(88, 8).Value = 28263.938
$ws.Cells.Item(88, 9).Value = 2044.125
$ws.Cells.Item(88, 10).Value = 54483.75
$ws.Cells.Item(88, 11).Value = 2044.125
$ws.Cells.Item(88, 12).Value = 54483.75
$ws.Cells.Item(88, 13).Value = -1638.125
$ws.Cells.Item(88, 14).Value = -55295.75
$ws.Cells.Item(89, 8).Value = 2191.875
$ws.Cells.Item(89, 9).Value = 2510.6
$ws.Cells.Item(89, 11).Value = 12553
$ws.Cells.Item(89, 13).Value = -6937
$ws.Cells.Item(90, 8).Value = 59570.645
$ws.Cells.Item(90, 10).Value = 59570.645
$ws.Cells.Item(90, 12).Value = 178711.935
$ws.Cells.Item(90, 14).Value = -191191.935
$ws.Cells.Item(91, 8).Value = 28263.938
$ws.Cells.Item(91, 9).Value = 2044.125
$ws.Cells.Item(91, 10).Value = 54483.75
$ws.Cells.Item(91, 11).Value = 2044.125
$ws.Cells.Item(91, 12).Value = 54483.75
$ws.Cells.Item(91, 13).Value = -640.125
$ws.Cells.Item(91, 14).Value = -57291.75
$ws.Cells.Item(98, 8).Value = 141584.08
$ws.Cells.Item(98, 9).Value = 116574.125
$ws.Cells.Item(98, 11).Value = 116574.125
$ws.Cells.Item(98, 13).Value = -115076.125
$ws.Cells.Item(106, 8).Value = 7043.8184
$ws.Cells.Item(106, 9).Value = 6532.8237
$ws.Cells.Item(106, 10).Value = 8781.2
$ws.Cells.Item(106, 11).Value = 6532.8237
$ws.Cells.Item(106, 12).Value = 8781.2
$ws.Cells.Item(106, 13).Value = -5901.8237
$ws.Cells.Item(106, 14).Value = -10043.2
$ws.Cells.Item(112, 8).Value = 2081.7036
$ws.Cells.Item(112, 10).Value = 2277.682
$ws.Cells.Item(112, 12).Value = 6833.045999999999
$ws.Cells.Item(112, 14).Value = -9049.045999999998
$ws.Cells.Item(116, 8).Value = 7142.3022
$ws.Cells.Item(116, 9).Value = 6749.2085
$ws.Cells.Item(116, 11).Value = 6749.2085
$ws.Cells.Item(116, 13).Value = -3307.2085
$ws.Cells.Item(122, 8).Value = 141584.08
$ws.Cells.Item(122, 9).Value = 116574.125
$ws.Cells.Item(122, 11).Value = 349722.375
$ws.Cells.Item(122, 13).Value = -347272.375
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 14).ClearContents()
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 14).ClearContents()
$ws.Cells.Item(131, 8).Value = 134556
$ws.Cells.Item(131, 9).Value = 134556
$ws.Cells.Item(131, 11).Value = 403668
$ws.Cells.Item(131, 13).Value = -398628
$ws.Cells.Item(132, 8).Value = 2010.0728
$ws.Cells.Item(132, 9).Value = 1524.5778
$ws.Cells.Item(132, 10).Value = 4194.8
$ws.Cells.Item(132, 11).Value = 4573.7334
$ws.Cells.Item(132, 12).Value = 12584.4
$ws.Cells.Item(132, 13).Value = -2043.7334
$ws.Cells.Item(132, 14).Value = -17644.4
$ws.Cells.Item(133, 8).Value = 180640
$ws.Cells.Item(133, 10).Value = 180640
$ws.Cells.Item(133, 12).Value = 180640
$ws.Cells.Item(133, 14).Value = -190760
$ws.Cells.Item(138, 8).Value = 3817.1633
$ws.Cells.Item(138, 9).Value = 1967.5454
$ws.Cells.Item(138, 10).Value = 4352.579
$ws.Cells.Item(138, 11).Value = 5902.6362
$ws.Cells.Item(138, 12).Value = 13057.737
$ws.Cells.Item(138, 13).Value = -762.6361999999999
$ws.Cells.Item(138, 14).Value = -23337.737
$ws.Cells.Item(139, 8).Value = 64017.445
$ws.Cells.Item(139, 10).Value = 64017.445
$ws.Cells.Item(139, 12).Value = 64017.445
$ws.Cells.Item(139, 14).Value = -74297.445
$ws.Cells.Item(140, 8).Value = 99542.29
$ws.Cells.Item(140, 10).Value = 99542.29
$ws.Cells.Item(140, 12).Value = 99542.29
$ws.Cells.Item(140, 14).Value = -109902.29

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6665.15
$ws.Cells.Item(32, 9).Value = 4409.4443
$ws.Cells.Item(32, 10).Value = 10505.946
$ws.Cells.Item(32, 11).Value = 4409.4443
$ws.Cells.Item(32, 12).Value = 10505.946
$ws.Cells.Item(32, 13).Value = -4122.4443
$ws.Cells.Item(32, 14).Value = -11079.946
$ws.Cells.Item(45, 8).Value = 1226.7858
$ws.Cells.Item(45, 9).Value = 1143.6364
$ws.Cells.Item(45, 10).Value = 1531.6666
$ws.Cells.Item(45, 11).Value = 1143.6364
$ws.Cells.Item(45, 12).Value = 1531.6666
$ws.Cells.Item(45, 13).Value = -766.6364000000001
$ws.Cells.Item(45, 14).Value = -2285.6666
$ws.Cells.Item(61, 8).Value = 26338.129
$ws.Cells.Item(61, 9).Value = 2582.8235
$ws.Cells.Item(61, 11).Value = 2582.8235
$ws.Cells.Item(61, 13).Value = -2370.8235
$ws.Cells.Item(88, 8).Value = 1081.3334
$ws.Cells.Item(88, 9).Value = 706.36365
$ws.Cells.Item(88, 11).Value = 706.36365
$ws.Cells.Item(88, 13).Value = -300.36365
$ws.Cells.Item(91, 8).Value = 1081.3334
$ws.Cells.Item(91, 9).Value = 706.36365
$ws.Cells.Item(91, 11).Value = 706.36365
$ws.Cells.Item(91, 13).Value = 697.63635
$ws.Cells.Item(102, 8).Value = 626789.94
$ws.Cells.Item(102, 9).Value = 911067.44
$ws.Cells.Item(102, 11).Value = 911067.44
$ws.Cells.Item(102, 13).Value = -909445.44
$ws.Cells.Item(110, 8).Value = 58097.375
$ws.Cells.Item(110, 9).Value = 68442.89
$ws.Cells.Item(110, 11).Value = 68442.89
$ws.Cells.Item(110, 13).Value = -66397.89
$ws.Cells.Item(122, 8).Value = 2235.1304
$ws.Cells.Item(122, 9).Value = 1850.9412
$ws.Cells.Item(122, 10).Value = 3323.6667
$ws.Cells.Item(122, 11).Value = 5552.8236
$ws.Cells.Item(122, 12).Value = 9971.000100000001
$ws.Cells.Item(122, 13).Value = -3102.8236
$ws.Cells.Item(122, 14).Value = -14871.0001
$ws.Cells.Item(132, 8).Value = 4065.6667
$ws.Cells.Item(132, 9).Value = 3588.9
$ws.Cells.Item(132, 11).Value = 10766.7
$ws.Cells.Item(132, 13).Value = -8236.7
$ws.Cells.Item(136, 8).Value = 26338.129
$ws.Cells.Item(136, 9).Value = 2582.8235
$ws.Cells.Item(136, 11).Value = 7748.470499999999
$ws.Cells.Item(136, 13).Value = -5198.470499999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(10, 8).Value = 1519.8
$ws.Cells.Item(10, 9).Value = 1100
$ws.Cells.Item(10, 10).Value = 1799.6666
$ws.Cells.Item(10, 11).Value = 1100
$ws.Cells.Item(10, 12).Value = 1799.6666
$ws.Cells.Item(10, 13).Value = -960
$ws.Cells.Item(10, 14).Value = -2079.6666
$ws.Cells.Item(12, 8).Value = 1175
$ws.Cells.Item(12, 9).Value = 400
$ws.Cells.Item(12, 10).Value = 3500
$ws.Cells.Item(12, 11).Value = 400
$ws.Cells.Item(12, 12).Value = 3500
$ws.Cells.Item(12, 13).Value = -232
$ws.Cells.Item(12, 14).Value = -3836
$ws.Cells.Item(86, 8).Value = 44207.668
$ws.Cells.Item(86, 9).Value = 1470.5555
$ws.Cells.Item(86, 10).Value = 76260.5
$ws.Cells.Item(86, 11).Value = 1470.5555
$ws.Cells.Item(86, 12).Value = 76260.5
$ws.Cells.Item(86, 13).Value = -347.5554999999999
$ws.Cells.Item(86, 14).Value = -78506.5
$ws.Cells.Item(89, 8).Value = 44207.668
$ws.Cells.Item(89, 9).Value = 1470.5555
$ws.Cells.Item(89, 10).Value = 76260.5
$ws.Cells.Item(89, 11).Value = 7352.7775
$ws.Cells.Item(89, 12).Value = 381302.5
$ws.Cells.Item(89, 13).Value = -1736.7775
$ws.Cells.Item(89, 14).Value = -392534.5
$ws.Cells.Item(94, 8).Value = 1147.3125
$ws.Cells.Item(94, 9).Value = 720.56
$ws.Cells.Item(94, 11).Value = 720.56
$ws.Cells.Item(94, 13).Value = -269.5599999999999
$ws.Cells.Item(105, 8).Value = 2514.84
$ws.Cells.Item(105, 9).Value = 2550.389
$ws.Cells.Item(105, 11).Value = 2550.389
$ws.Cells.Item(105, 13).Value = -803.3890000000001
$ws.Cells.Item(134, 8).Value = 30884
$ws.Cells.Item(134, 9).Value = 33362.484
$ws.Cells.Item(134, 10).Value = 6925.3335
$ws.Cells.Item(134, 11).Value = 100087.452
$ws.Cells.Item(134, 12).Value = 20776.0005
$ws.Cells.Item(134, 13).Value = -97552.45199999999
$ws.Cells.Item(134, 14).Value = -25846.0005
$ws.Cells.Item(138, 8).Value = 86853.05
$ws.Cells.Item(138, 10).Value = 86853.05
$ws.Cells.Item(138, 12).Value = 86853.05
$ws.Cells.Item(138, 14).Value = -97133.05

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1590
$ws.Cells.Item(16, 9).Value = 858.7
$ws.Cells.Item(16, 10).Value = 3052.6
$ws.Cells.Item(16, 11).Value = 858.7
$ws.Cells.Item(16, 12).Value = 3052.6
$ws.Cells.Item(16, 13).Value = -571.7
$ws.Cells.Item(16, 14).Value = -3626.6
$ws.Cells.Item(31, 8).Value = 2984.6965
$ws.Cells.Item(31, 9).Value = 2100.7693
$ws.Cells.Item(31, 11).Value = 2100.7693
$ws.Cells.Item(31, 13).Value = -1805.7693
$ws.Cells.Item(34, 8).Value = 2984.6965
$ws.Cells.Item(34, 9).Value = 2100.7693
$ws.Cells.Item(34, 11).Value = 2100.7693
$ws.Cells.Item(34, 13).Value = -1898.7693
$ws.Cells.Item(99, 8).Value = 7550.522
$ws.Cells.Item(99, 9).Value = 6215.7646
$ws.Cells.Item(99, 11).Value = 6215.7646
$ws.Cells.Item(99, 13).Value = -4717.7646
$ws.Cells.Item(105, 8).Value = 806.4091
$ws.Cells.Item(105, 9).Value = 802.2778
$ws.Cells.Item(105, 10).Value = 825
$ws.Cells.Item(105, 11).Value = 802.2778
$ws.Cells.Item(105, 12).Value = 825
$ws.Cells.Item(105, 13).Value = 944.7222
$ws.Cells.Item(105, 14).Value = -4319
$ws.Cells.Item(113, 8).Value = 1590
$ws.Cells.Item(113, 9).Value = 858.7
$ws.Cells.Item(113, 10).Value = 3052.6
$ws.Cells.Item(113, 11).Value = 858.7
$ws.Cells.Item(113, 12).Value = 3052.6
$ws.Cells.Item(113, 13).Value = 1311.3
$ws.Cells.Item(113, 14).Value = -7392.6
$ws.Cells.Item(122, 8).Value = 1839.4642
$ws.Cells.Item(122, 9).Value = 1535.5454
$ws.Cells.Item(122, 11).Value = 4606.6362
$ws.Cells.Item(122, 13).Value = -2156.6362
$ws.Cells.Item(126, 8).Value = 7550.522
$ws.Cells.Item(126, 9).Value = 6215.7646
$ws.Cells.Item(126, 11).Value = 18647.2938
$ws.Cells.Item(126, 13).Value = -16177.2938
$ws.Cells.Item(134, 8).Value = 1595.45
$ws.Cells.Item(134, 9).Value = 874.7292
$ws.Cells.Item(134, 11).Value = 2624.1876
$ws.Cells.Item(134, 13).Value = -89.1876000000002
$ws.Cells.Item(135, 8).Value = 66933.336
$ws.Cells.Item(135, 10).Value = 66933.336
$ws.Cells.Item(135, 12).Value = 66933.336
$ws.Cells.Item(135, 14).Value = -77073.336
$ws.Cells.Item(138, 8).Value = 92254.14
$ws.Cells.Item(138, 10).Value = 92254.14
$ws.Cells.Item(138, 12).Value = 92254.14
$ws.Cells.Item(138, 14).Value = -102534.14

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 749.5
$ws.Cells.Item(8, 9).Value = 749.5
$ws.Cells.Item(8, 11).Value = 2248.5
$ws.Cells.Item(8, 13).Value = -2109.5
$ws.Cells.Item(12, 8).Value = 827.0789
$ws.Cells.Item(12, 10).Value = 1030.6207
$ws.Cells.Item(12, 12).Value = 3091.8621
$ws.Cells.Item(12, 14).Value = -3437.8621
$ws.Cells.Item(14, 8).Value = 2425
$ws.Cells.Item(14, 9).Value = 2425
$ws.Cells.Item(14, 11).Value = 7275
$ws.Cells.Item(14, 13).Value = -7102
$ws.Cells.Item(23, 8).Value = 676.9167
$ws.Cells.Item(23, 10).Value = 653.8333
$ws.Cells.Item(23, 12).Value = 1961.4999
$ws.Cells.Item(23, 14).Value = -2431.4999
$ws.Cells.Item(34, 8).Value = 298
$ws.Cells.Item(34, 10).Value = 297
$ws.Cells.Item(34, 12).Value = 891
$ws.Cells.Item(34, 14).Value = -1059
$ws.Cells.Item(37, 8).Value = 120861.07
$ws.Cells.Item(37, 10).Value = 120861.07
$ws.Cells.Item(37, 12).Value = 362583.21
$ws.Cells.Item(37, 14).Value = -362807.21
$ws.Cells.Item(39, 8).Value = 12952.846
$ws.Cells.Item(39, 9).Value = 1999.6666
$ws.Cells.Item(39, 10).Value = 16238.8
$ws.Cells.Item(39, 11).Value = 5998.9998
$ws.Cells.Item(39, 12).Value = 48716.39999999999
$ws.Cells.Item(39, 13).Value = -5704.9998
$ws.Cells.Item(39, 14).Value = -49304.39999999999
$ws.Cells.Item(44, 8).Value = 901.0909
$ws.Cells.Item(44, 9).Value = 546.44446
$ws.Cells.Item(44, 10).Value = 2497
$ws.Cells.Item(44, 11).Value = 1639.33338
$ws.Cells.Item(44, 12).Value = 7491
$ws.Cells.Item(44, 13).Value = -1241.33338
$ws.Cells.Item(44, 14).Value = -8287
$ws.Cells.Item(55, 8).Value = 7553.625
$ws.Cells.Item(55, 9).Value = 700
$ws.Cells.Item(55, 10).Value = 11665.8
$ws.Cells.Item(55, 11).Value = 2100
$ws.Cells.Item(55, 12).Value = 34997.39999999999
$ws.Cells.Item(55, 13).Value = -1923
$ws.Cells.Item(55, 14).Value = -35351.39999999999
$ws.Cells.Item(97, 8).Value = 2975
$ws.Cells.Item(97, 10).Value = 5500
$ws.Cells.Item(97, 12).Value = 16500
$ws.Cells.Item(97, 14).Value = -17492
$ws.Cells.Item(98, 8).Value = 975.9
$ws.Cells.Item(98, 9).Value = 999.5
$ws.Cells.Item(98, 10).Value = 960.1667
$ws.Cells.Item(98, 11).Value = 2998.5
$ws.Cells.Item(98, 12).Value = 2880.5001
$ws.Cells.Item(98, 13).Value = -1500.5
$ws.Cells.Item(98, 14).Value = -5876.5001
$ws.Cells.Item(113, 8).Value = 71344.62
$ws.Cells.Item(113, 10).Value = 102620.22
$ws.Cells.Item(113, 12).Value = 307860.66
$ws.Cells.Item(113, 14).Value = -312200.66
$ws.Cells.Item(131, 8).Value = 1833.4375
$ws.Cells.Item(131, 9).Value = 1717.9
$ws.Cells.Item(131, 10).Value = 1885.9546
$ws.Cells.Item(131, 11).Value = 5153.700000000001
$ws.Cells.Item(131, 12).Value = 5657.8638
$ws.Cells.Item(131, 13).Value = -113.7000000000007
$ws.Cells.Item(131, 14).Value = -15737.8638
$ws.Cells.Item(132, 8).Value = 2940.2778
$ws.Cells.Item(132, 9).Value = 1512.7142
$ws.Cells.Item(132, 11).Value = 13614.4278
$ws.Cells.Item(132, 13).Value = -11084.4278
$ws.Cells.Item(137, 8).Value = 110128.6
$ws.Cells.Item(137, 10).Value = 4675.6196
$ws.Cells.Item(137, 12).Value = 14026.8588
$ws.Cells.Item(137, 14).Value = -24226.8588

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 8652.655
$ws.Cells.Item(80, 9).Value = 3675.6667
$ws.Cells.Item(80, 11).Value = 3675.6667
$ws.Cells.Item(80, 13).Value = -2677.6667
$ws.Cells.Item(83, 8).Value = 8652.655
$ws.Cells.Item(83, 9).Value = 3675.6667
$ws.Cells.Item(83, 11).Value = 18378.3335
$ws.Cells.Item(83, 13).Value = -13386.3335
$ws.Cells.Item(97, 8).Value = 357.77777
$ws.Cells.Item(97, 9).Value = 182.5
$ws.Cells.Item(97, 10).Value = 708.3333
$ws.Cells.Item(97, 11).Value = 182.5
$ws.Cells.Item(97, 12).Value = 708.3333
$ws.Cells.Item(97, 13).Value = 313.5
$ws.Cells.Item(97, 14).Value = -1700.3333
$ws.Cells.Item(107, 8).Value = 221.25581
$ws.Cells.Item(107, 9).Value = 233.6923
$ws.Cells.Item(107, 10).Value = 202.23529
$ws.Cells.Item(107, 11).Value = 233.6923
$ws.Cells.Item(107, 12).Value = 202.23529
$ws.Cells.Item(107, 13).Value = 1686.3077
$ws.Cells.Item(107, 14).Value = -4042.23529
$ws.Cells.Item(113, 8).Value = 9919.077
$ws.Cells.Item(113, 9).Value = 11268.454
$ws.Cells.Item(113, 10).Value = 2497.5
$ws.Cells.Item(113, 11).Value = 11268.454
$ws.Cells.Item(113, 12).Value = 2497.5
$ws.Cells.Item(113, 13).Value = -9098.454
$ws.Cells.Item(113, 14).Value = -6837.5
$ws.Cells.Item(132, 8).Value = 2250.2144
$ws.Cells.Item(132, 9).Value = 1950.0857
$ws.Cells.Item(132, 10).Value = 3750.8572
$ws.Cells.Item(132, 11).Value = 5850.257100000001
$ws.Cells.Item(132, 12).Value = 11252.5716
$ws.Cells.Item(132, 13).Value = -3320.257100000001
$ws.Cells.Item(132, 14).Value = -16312.5716
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(135, 8).Value = 99658
$ws.Cells.Item(135, 10).Value = 99658
$ws.Cells.Item(135, 12).Value = 99658
$ws.Cells.Item(135, 14).Value = -109798

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 45131.57
$ws.Cells.Item(7, 9).Value = 62283.4
$ws.Cells.Item(7, 11).Value = 62283.4
$ws.Cells.Item(7, 13).Value = -62171.4
$ws.Cells.Item(40, 8).Value = 30882.41
$ws.Cells.Item(40, 9).Value = 35205.742
$ws.Cells.Item(40, 10).Value = 14165.533
$ws.Cells.Item(40, 11).Value = 35205.742
$ws.Cells.Item(40, 12).Value = 14165.533
$ws.Cells.Item(40, 13).Value = -35069.742
$ws.Cells.Item(40, 14).Value = -14437.533
$ws.Cells.Item(55, 8).Value = 690.3214
$ws.Cells.Item(55, 9).Value = 561
$ws.Cells.Item(55, 11).Value = 561
$ws.Cells.Item(55, 13).Value = -388
$ws.Cells.Item(61, 8).Value = 5032.5713
$ws.Cells.Item(61, 9).Value = 5097.2
$ws.Cells.Item(61, 11).Value = 5097.2
$ws.Cells.Item(61, 13).Value = -4895.2
$ws.Cells.Item(68, 8).Value = 2929.7778
$ws.Cells.Item(68, 9).Value = 2222
$ws.Cells.Item(68, 10).Value = 4770
$ws.Cells.Item(68, 11).Value = 2222
$ws.Cells.Item(68, 12).Value = 4770
$ws.Cells.Item(68, 13).Value = -1473
$ws.Cells.Item(68, 14).Value = -6268
$ws.Cells.Item(71, 8).Value = 2929.7778
$ws.Cells.Item(71, 9).Value = 2222
$ws.Cells.Item(71, 10).Value = 4770
$ws.Cells.Item(71, 11).Value = 11110
$ws.Cells.Item(71, 12).Value = 23850
$ws.Cells.Item(71, 13).Value = -7366
$ws.Cells.Item(71, 14).Value = -31338
$ws.Cells.Item(93, 8).Value = 1762.091
$ws.Cells.Item(93, 9).Value = 1244.5
$ws.Cells.Item(93, 10).Value = 3142.3333
$ws.Cells.Item(93, 11).Value = 1244.5
$ws.Cells.Item(93, 12).Value = 3142.3333
$ws.Cells.Item(93, 13).Value = 3.5
$ws.Cells.Item(93, 14).Value = -5638.3333
$ws.Cells.Item(113, 8).Value = 5032.5713
$ws.Cells.Item(113, 9).Value = 5097.2
$ws.Cells.Item(113, 11).Value = 5097.2
$ws.Cells.Item(113, 13).Value = -2927.2
$ws.Cells.Item(122, 8).Value = 3380.4517
$ws.Cells.Item(122, 9).Value = 3352.8965
$ws.Cells.Item(122, 10).Value = 3780
$ws.Cells.Item(122, 11).Value = 10058.6895
$ws.Cells.Item(122, 12).Value = 11340
$ws.Cells.Item(122, 13).Value = -7608.6895
$ws.Cells.Item(122, 14).Value = -16240
$ws.Cells.Item(126, 8).Value = 45131.57
$ws.Cells.Item(126, 9).Value = 62283.4
$ws.Cells.Item(126, 11).Value = 186850.2
$ws.Cells.Item(126, 13).Value = -184380.2
$ws.Cells.Item(135, 8).Value = 77863.75
$ws.Cells.Item(135, 10).Value = 77863.75
$ws.Cells.Item(135, 12).Value = 77863.75
$ws.Cells.Item(135, 14).Value = -88003.75
$ws.Cells.Item(136, 8).Value = 3935.7446
$ws.Cells.Item(136, 9).Value = 3120.054
$ws.Cells.Item(136, 11).Value = 9360.162
$ws.Cells.Item(136, 13).Value = -6810.162

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 159091.47
$ws.Cells.Item(62, 9).Value = 171990.75
$ws.Cells.Item(62, 11).Value = 171990.75
$ws.Cells.Item(62, 13).Value = -171366.75
$ws.Cells.Item(65, 8).Value = 159091.47
$ws.Cells.Item(65, 9).Value = 171990.75
$ws.Cells.Item(65, 11).Value = 859953.75
$ws.Cells.Item(65, 13).Value = -856833.75
$ws.Cells.Item(70, 8).Value = 25081.4
$ws.Cells.Item(70, 10).Value = 26938.363
$ws.Cells.Item(70, 14).Value = -27568.363
$ws.Cells.Item(73, 8).Value = 25081.4
$ws.Cells.Item(73, 10).Value = 26938.363
$ws.Cells.Item(73, 14).Value = -29122.363
$ws.Cells.Item(81, 8).Value = 1327.6666
$ws.Cells.Item(81, 9).Value = 1374.875
$ws.Cells.Item(81, 11).Value = 2749.75
$ws.Cells.Item(81, 13).Value = -1688.75
$ws.Cells.Item(84, 8).Value = 1327.6666
$ws.Cells.Item(84, 9).Value = 1374.875
$ws.Cells.Item(84, 11).Value = 13748.75
$ws.Cells.Item(84, 13).Value = -8444.75
$ws.Cells.Item(96, 8).Value = 2000.35
$ws.Cells.Item(96, 9).Value = 1994.7693
$ws.Cells.Item(96, 11).Value = 1994.7693
$ws.Cells.Item(96, 13).Value = -621.7692999999999
$ws.Cells.Item(113, 8).Value = 1961.1578
$ws.Cells.Item(113, 9).Value = 1921
$ws.Cells.Item(113, 10).Value = 2005.7778
$ws.Cells.Item(113, 11).Value = 5763
$ws.Cells.Item(113, 12).Value = 6017.3334
$ws.Cells.Item(113, 13).Value = -3593
$ws.Cells.Item(113, 14).Value = -10357.3334
$ws.Cells.Item(122, 8).Value = 2245.7837
$ws.Cells.Item(122, 9).Value = 2109.7273
$ws.Cells.Item(122, 11).Value = 6329.1819
$ws.Cells.Item(122, 13).Value = -3879.1819
$ws.Cells.Item(132, 8).Value = 7469.9614
$ws.Cells.Item(132, 9).Value = 7283.2
$ws.Cells.Item(132, 10).Value = 7724.636
$ws.Cells.Item(132, 11).Value = 21849.6
$ws.Cells.Item(132, 12).Value = 23173.908
$ws.Cells.Item(132, 13).Value = -19319.6
$ws.Cells.Item(132, 14).Value = -28233.908
$ws.Cells.Item(136, 8).Value = 2203.9033
$ws.Cells.Item(136, 9).Value = 2235.9434
$ws.Cells.Item(136, 11).Value = 6707.8302
$ws.Cells.Item(136, 13).Value = -4157.8302
$ws.Cells.Item(139, 8).Value = 87273.5
$ws.Cells.Item(139, 9).Value = 49925
$ws.Cells.Item(139, 10).Value = 97944.5
$ws.Cells.Item(139, 11).Value = 49925
$ws.Cells.Item(139, 12).Value = 97944.5
$ws.Cells.Item(139, 13).Value = -44785
$ws.Cells.Item(139, 14).Value = -108224.5
$ws.Cells.Item(141, 8).Value = 96985.336
$ws.Cells.Item(141, 10).Value = 96985.336
$ws.Cells.Item(141, 12).Value = 96985.336
$ws.Cells.Item(141, 14).Value = -107345.336
